# Fixed biometrics import error
# Several biometrics time-string values imported onto the DTR Summary sheet
# were missing a trailing "0" (e.g. "0.4.3" should really be "0.4.30").
# Every cell holding one of the affected values needs to be corrected so the
# workbook no longer shares the old (wrong) text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of old (incorrect) biometrics values to the corrected values, together
# with every cell address in the sheet that currently holds that value.
$fixes = @(
    @{ Old = "0.4.3"; New = "0.4.30"; Cells = @("E4", "F44") },
    @{ Old = "0.0.3"; New = "0.0.30"; Cells = @("E10", "E19", "F29", "F48", "F55") },
    @{ Old = "0.1.3"; New = "0.1.30"; Cells = @("E24", "F25", "F46") },
    @{ Old = "0.3.3"; New = "0.3.30"; Cells = @("E49") },
    @{ Old = "0.2.3"; New = "0.2.30"; Cells = @("E50") }
)

foreach ($fix in $fixes) {
    foreach ($addr in $fix.Cells) {
        $cell = $ws.Range($addr)
        $current = [string]$cell.Value2
        if ($current -eq $fix.Old) {
            $cell.Value2 = $fix.New
        }
    }
}
